# Apply the "Add files via upload" edit to ff_pub.xlsx:
#  - rename the sheet from "Sheet1" to "ff_pub"
#  - append two new movie rows (135 & 136) to the tracking table
#  - leave the selection positioned over the newly-added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "ff_pub"

# Row 135 currently exists only as an empty, custom-formatted placeholder
# row (s="3" customFormat="1") -- clear that formatting before writing data
# so the saved row looks like a normal data row (matching rows above it).
$ws.Rows.Item(135).ClearFormats()

$ws.Range("A135").Value = "Cross of Iron (1977)"
$ws.Range("B135").Value = 4.5
$ws.Range("C135").Value = 4
$ws.Range("D135").Value = 4.5
$ws.Range("E135").Value = "Crosses of Irons"
$ws.Range("F135").Value = "https://www.reddit.com/r/FriendlyFirePodcast/comments/g75kk9/ep_120_cross_of_iron_1977/"
$ws.Range("H135").Value = 120
$ws.Range("J135").Value = "https://www.justwatch.com/us/movie/cross-of-iron"
$ws.Range("L135").Value = "It is 1943, and the German army—ravaged and demoralised—is hastily retreating from the Russian front. In the midst of the madness, conflict brews between the aristocratic yet ultimately pusillanimous Captain Stransky and the courageous Corporal Steiner. Stransky is the only man who believes that the Third Reich is still vastly superior to the Russian army. However, within his pompous persona lies a quivering coward who longs for the Iron Cross so that he can return to Berlin a hero. Steiner, on the other hand is cynical, defiantly non-conformist and more concerned with the safety of his own men rather than the horde of military decorations offered to him by his superiors."
$ws.Range("M135").Value = "https://www.imdb.com/title/tt0074695/"
$ws.Range("N135").Value = 1977

# Row 136 is a brand-new row.
$ws.Range("A136").Value = "Good Morning, Vietnam (1987)"
$ws.Range("B136").Value = 4.5
$ws.Range("C136").Value = 3.5
$ws.Range("D136").Value = 3
$ws.Range("E136").Value = "Green suits"
$ws.Range("F136").Value = "https://www.reddit.com/r/FriendlyFirePodcast/comments/gbf0bp/ep_121_good_morning_vietnam_1987/"
$ws.Range("H136").Value = 121
$ws.Range("J136").Value = "https://www.justwatch.com/us/movie/good-morning-vietnam"
$ws.Range("L136").Value = "In 1965, an unorthodox and irreverent DJ named Adrian Cronauer begins to shake up things when he is assigned to the U.S. Armed Services radio station in Vietnam."
$ws.Range("M136").Value = "https://www.imdb.com/title/tt0093105/"
$ws.Range("N136").Value = 1987

# Select the newly added rows, mirroring the author's final selection state.
$ws.Range("A135:XFD136").Select()
